# Updated cryptos list on Fri Mar 24 21:57:14 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) text values for the crypto table, and
# swaps the NEARProtocol / Quant rows (48 & 49) which changed rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # The sheet stores Price/Link/Coin columns as plain text (inline strings),
    # even though several values look numeric (e.g. "1.003"). Force the
    # NumberFormat to Text before assigning so Excel does not silently
    # reinterpret the string as a number, then restore the default style so
    # no stray formatting is left behind on the cell.
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.Style = 'Normal'
}

Set-TextCell 'D2' '27.399.88'
Set-TextCell 'E2' '  -2.98%  '
Set-TextCell 'D3' '1.748.75'
Set-TextCell 'E3' '  -3.83%  '
Set-TextCell 'D4' '1.003'
Set-TextCell 'E4' '  +0.07%  '
Set-TextCell 'D5' '321.21'
Set-TextCell 'E5' '  -2.60%  '
Set-TextCell 'D6' '1.001'
Set-TextCell 'E6' '  +0.13%  '
Set-TextCell 'D7' '0.4225'
Set-TextCell 'E7' '  -4.96%  '
Set-TextCell 'D8' '0.3592'
Set-TextCell 'E8' '  -2.96%  '
Set-TextCell 'D9' '0.07526'
Set-TextCell 'E9' '  -2.25%  '
Set-TextCell 'D10' '42.36'
Set-TextCell 'E10' '  -5.02%  '
Set-TextCell 'D11' '1.091'
Set-TextCell 'E11' '  -3.17%  '
Set-TextCell 'D12' '1.002'
Set-TextCell 'E12' '  +0.13%  '
Set-TextCell 'D13' '20.64'
Set-TextCell 'E13' '  -6.73%  '
Set-TextCell 'D14' '6.029'
Set-TextCell 'E14' '  -3.90%  '
Set-TextCell 'D15' '7.201'
Set-TextCell 'E15' '  -5.07%  '
Set-TextCell 'D16' '1.749.18'
Set-TextCell 'E16' '  -5.45%  '
Set-TextCell 'D17' '91.33'
Set-TextCell 'E17' '  -1.46%  '
Set-TextCell 'D18' '0.00001067'
Set-TextCell 'E18' '  -1.62%  '
Set-TextCell 'D19' '0.06347'
Set-TextCell 'E19' '  -3.93%  '
Set-TextCell 'D20' '1.003'
Set-TextCell 'E20' '  +0.21%  '
Set-TextCell 'D21' '17.02'
Set-TextCell 'E21' '  -3.05%  '
Set-TextCell 'D22' '5.876'
Set-TextCell 'E22' '  -5.23%  '
Set-TextCell 'D23' '27.444.11'
Set-TextCell 'E23' '  -3.04%  '
Set-TextCell 'D24' '11.17'
Set-TextCell 'E24' '  -4.43%  '
Set-TextCell 'D25' '2.090'
Set-TextCell 'E25' '  -1.86%  '
Set-TextCell 'D26' '161.10'
Set-TextCell 'E26' '  +3.26%  '
Set-TextCell 'D27' '20.24'
Set-TextCell 'E27' '  -2.73%  '
Set-TextCell 'D28' '1.953.72'
Set-TextCell 'E28' '  -4.30%  '
Set-TextCell 'D29' '2.126'
Set-TextCell 'E29' '  -8.54%  '
Set-TextCell 'D30' '123.04'
Set-TextCell 'E30' '  -4.11%  '
Set-TextCell 'D31' '1.107'
Set-TextCell 'E31' '  -8.30%  '
Set-TextCell 'D32' '3.642'
Set-TextCell 'E32' '  -0.07%  '
Set-TextCell 'D33' '5.536'
Set-TextCell 'E33' '  -5.64%  '
Set-TextCell 'D34' '0.08829'
Set-TextCell 'E34' '  -4.33%  '
Set-TextCell 'E35' '  -6.93%  '
Set-TextCell 'E36' '  -3.45%  '
Set-TextCell 'D37' '0.2097'
Set-TextCell 'E37' '  -3.92%  '
Set-TextCell 'D38' '0.05993'
Set-TextCell 'E38' '  -3.70%  '
Set-TextCell 'D39' '0.6315'
Set-TextCell 'E39' '  -4.02%  '
Set-TextCell 'D40' '4.924'
Set-TextCell 'E40' '  -4.43%  '
Set-TextCell 'D41' '1.180'
Set-TextCell 'E41' '  -1.63%  '
Set-TextCell 'D42' '1.001'
Set-TextCell 'E42' '  +0.14%  '
Set-TextCell 'D43' '7.852'
Set-TextCell 'E43' '  -3.94%  '
Set-TextCell 'D44' '1.388'
Set-TextCell 'E44' '  -0.20%  '
Set-TextCell 'D45' '13.40'
Set-TextCell 'E45' '  -3.57%  '
Set-TextCell 'D46' '0.5846'
Set-TextCell 'E46' '  -3.80%  '
Set-TextCell 'D47' '3.686'
Set-TextCell 'E47' '  -2.23%  '
Set-TextCell 'B48' 'Quant'
Set-TextCell 'C48' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 'D48' '122.33'
Set-TextCell 'E48' '  -3.95%  '
Set-TextCell 'B49' 'NEARProtocol'
Set-TextCell 'C49' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D49' '1.966'
Set-TextCell 'E49' '  -3.63%  '
Set-TextCell 'D50' '1.151'
Set-TextCell 'E50' '  -0.25%  '
Set-TextCell 'D51' '0.06801'
Set-TextCell 'E51' '  -2.58%  '
